$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1, styled like the other header cells (bold, centered, bordered)
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data cells F2:F36 - plain text timestamps
$ws.Range("F2").Value = "2021-10-05 10:52:38.702113"
$ws.Range("F3").Value = "2021-10-05 10:52:38.702124"
$ws.Range("F4").Value = "2021-10-05 10:52:38.702128"
$ws.Range("F5").Value = "2021-10-05 10:52:38.702130"
$ws.Range("F6").Value = "2021-10-05 10:52:38.702133"
$ws.Range("F7").Value = "2021-10-05 10:52:38.702136"
$ws.Range("F8").Value = "2021-10-05 10:52:38.702138"
$ws.Range("F9").Value = "2021-10-05 10:52:38.702141"
$ws.Range("F10").Value = "2021-10-05 10:52:38.702144"
$ws.Range("F11").Value = "2021-10-05 10:52:38.702146"
$ws.Range("F12").Value = "2021-10-05 10:52:38.702149"
$ws.Range("F13").Value = "2021-10-05 10:52:38.702151"
$ws.Range("F14").Value = "2021-10-05 10:52:38.702154"
$ws.Range("F15").Value = "2021-10-05 10:52:38.702156"
$ws.Range("F16").Value = "2021-10-05 10:52:38.702159"
$ws.Range("F17").Value = "2021-10-05 10:52:38.702161"
$ws.Range("F18").Value = "2021-10-05 10:52:38.702164"
$ws.Range("F19").Value = "2021-10-05 10:52:38.702167"
$ws.Range("F20").Value = "2021-10-05 10:52:38.702170"
$ws.Range("F21").Value = "2021-10-05 10:52:38.702172"
$ws.Range("F22").Value = "2021-10-05 10:52:38.702175"
$ws.Range("F23").Value = "2021-10-05 10:52:38.702178"
$ws.Range("F24").Value = "2021-10-05 10:52:38.702180"
$ws.Range("F25").Value = "2021-10-05 10:52:38.702183"
$ws.Range("F26").Value = "2021-10-05 10:52:38.702186"
$ws.Range("F27").Value = "2021-10-05 10:52:38.702188"
$ws.Range("F28").Value = "2021-10-05 10:52:38.702191"
$ws.Range("F29").Value = "2021-10-05 10:52:38.702194"
$ws.Range("F30").Value = "2021-10-05 10:52:38.702196"
$ws.Range("F31").Value = "2021-10-05 10:52:38.702199"
$ws.Range("F32").Value = "2021-10-05 10:52:38.702201"
$ws.Range("F33").Value = "2021-10-05 10:52:38.702204"
$ws.Range("F34").Value = "2021-10-05 10:52:38.702207"
$ws.Range("F35").Value = "2021-10-05 10:52:38.702210"
$ws.Range("F36").Value = "2021-10-05 10:52:38.702212"
